# Updates the cryptos list figures (price/volume) to the latest scrape.
# Each (cell -> new text) pair below is applied via Range.Value. Cells whose
# new text parses as a plain number (e.g. "41.00", "6.01") are written with
# a leading "'" so Excel keeps storing them as text (matching the inlineStr
# type used throughout this sheet) instead of silently converting them to
# numbers and losing the exact formatting (trailing zeros, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.122.92"
$ws.Range("E2").Value = "  +5.84%  "
$ws.Range("D3").Value = "2.243.78"
$ws.Range("E3").Value = "  +5.25%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.54"
$ws.Range("E5").Value = "  +6.65%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +3.78%  "
$ws.Range("D7").Value = "'75.38"
$ws.Range("E7").Value = "  +10.77%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'0.604"
$ws.Range("E9").Value = "  +9.06%  "
$ws.Range("D10").Value = "'41.00"
$ws.Range("E10").Value = "  +9.62%  "
$ws.Range("D11").Value = "'0.0932"
$ws.Range("E11").Value = "  +5.79%  "
$ws.Range("D12").Value = "'6.91"
$ws.Range("E12").Value = "  +7.04%  "
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "2.578.85"
$ws.Range("E14").Value = "  +5.21%  "
$ws.Range("D15").Value = "'14.62"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "2.250.35"
$ws.Range("E16").Value = "  +5.73%  "
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = "  +3.48%  "
$ws.Range("D18").Value = "42.997.57"
$ws.Range("E18").Value = "  +5.96%  "
$ws.Range("E19").Value = "  +8.08%  "
$ws.Range("D20").Value = "'71.25"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("D21").Value = "'6.01"
$ws.Range("E21").Value = "  +6.37%  "
$ws.Range("E22").Value = "  +5.32%  "
$ws.Range("E23").Value = "  +20.56%  "
$ws.Range("D24").Value = "'229.93"
$ws.Range("E24").Value = "  +4.25%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'10.83"
$ws.Range("E26").Value = "  +4.64%  "
$ws.Range("D27").Value = "'3.43"
$ws.Range("E27").Value = "  +8.19%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'39.32"
$ws.Range("E29").Value = "  +32.00%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +4.15%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'171.66"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("D32").Value = "'20.25"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  +9.16%  "
$ws.Range("D34").Value = "'5.29"
$ws.Range("E34").Value = "  +7.49%  "
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  +14.86%  "
$ws.Range("D37").Value = "'4.47"
$ws.Range("E37").Value = "  +12.57%  "
$ws.Range("D38").Value = "'0.0329"
$ws.Range("E38").Value = "  +21.61%  "
$ws.Range("D39").Value = "'13.09"
$ws.Range("E39").Value = "  +17.33%  "
$ws.Range("E40").Value = "  +5.96%  "
$ws.Range("E41").Value = "  +13.66%  "
$ws.Range("D42").Value = "'5.42"
$ws.Range("E42").Value = "  +4.53%  "
$ws.Range("D43").Value = "'59.44"
$ws.Range("E43").Value = "  +6.38%  "
$ws.Range("D44").Value = "'104.68"
$ws.Range("E44").Value = "  +10.75%  "
$ws.Range("D45").Value = "'8.69"
$ws.Range("E45").Value = "  +7.99%  "
$ws.Range("D46").Value = "'0.482"
$ws.Range("E46").Value = "  +37.55%  "
$ws.Range("D47").Value = "'0.0993"
$ws.Range("E47").Value = "  +5.91%  "
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = "  +15.47%  "
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  +5.56%  "
$ws.Range("E50").Value = "  +6.54%  "
$ws.Range("E51").Value = "  +3.79%  "

Write-Output "Updated 96 cells (33 forced to text)."
